$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDescription = "控制元件各個部分的顏色`nbackgroundColor: 整個元件的背景顏色，包括日期按鈕和切換年份視窗的背景`nheaderColor: header 的背景顏色`nheaderTextColor: header 的文字顏色和切換月份的 icon 顏色`nchangeYearModalColor: 切換年份視窗的主色`nweekDaysColor: 星期文字的顏色`ndateTextColor: 日期按鈕中的文字顏色 (6碼HEX)`nselectedDateColor: 當日期按鈕被active 時的文字顏色 (6碼HEX)`nselectedDateBackgroundColor: 當日期按鈕被 active 時的背景顏色 (6碼HEX)`nconfirmButtonColor: 確認按鈕的文字顏色"

$ws.Range("C14").Value = $newDescription

$ws.Range("C17").Select()
